# Update indicator metadata on the "Пример" worksheet:
#  - B4:  indicator name text revised (age range 36-59 months made explicit)
#  - B6:  responsible organization's sub-division renamed
#  - B7:  contact person changed
#  - B8:  contact e-mail changed
#  - B9:  contact phone changed
#  - B10: organization website changed

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Пример")

$ws.Range("B4").Value = "4.2.1. Доля детей в возрасте от 36 до 59 месяцев, развивающихся без отклонений в плане здоровья, обучения и психосоциального благополучия, в разбивке по полу"
$ws.Range("B6").Value = "Национальный статистический комитет Кыргызской Республики (Управление статистики домашних хозяйств)"
$ws.Range("B7").Value = "Калымбетова Ы.И."
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

$ws.Range("B6").Select()
